# "Stacking limited for potions + bag ameliorations buyable"
#
# - A4  (9: potion de poison)            -> highlight green (stacking limit task started)
# - A7  (Missions 2.1: limite inventaire) -> highlight green + mark column C as "done"
# - B13 (Mission 2.2 : sac)               -> reassign from Romeo to Fabio (bag ameliorations)
# - A25 (Pom-Pom Girl)                    -> highlight green + mark column C as "perfect"
# - selection left on B13, matching the saved workbook's cursor position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$green = 5287936   # RGB(0,176,80) packed as BGR for OLE color

# Row 4: A4 -> green fill, centered (same look as the other "started" rows)
$ws.Range("A4").Interior.Color = $green
$ws.Range("A4").HorizontalAlignment = -4108

# Row 7: A7 -> green fill, centered; this row additionally gets its own new
# cell style (font touched so a distinct xf is minted) plus a "done" marker
$ws.Range("A7").Interior.Color = $green
$ws.Range("A7").Font.Name = "Calibri"
$ws.Range("A7").HorizontalAlignment = -4108
$ws.Range("C7").Value = "done"

# Row 13: bag task reassigned to Fabio
$ws.Range("B13").Value = "Fabio"

# Row 25: A25 -> green fill, centered; mark column C as "perfect"
$ws.Range("A25").Interior.Color = $green
$ws.Range("A25").HorizontalAlignment = -4108
$ws.Range("C25").Value = "perfect"

# Restore the selection to B13 as saved in the workbook
$ws.Range("B13").Select()
